$d = $word.ActiveDocument

# 1) Header cell: "Characteristic" -> "Baseline Characteristics", and un-bold it.
#    Using Find/Replacement.Font so only the matched run's formatting changes
#    (avoids also stamping the paragraph mark).
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Replacement.ClearFormatting()
$find1.Text = "Characteristic"
$find1.Replacement.Text = "Baseline Characteristics"
$find1.Replacement.Font.Bold = 0
$find1.Execute($find1.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 1, $false, $false, $false, $false) | Out-Null

# 2) Body cell: "Nonwhite" -> "Non-white" (formatting unchanged)
$d.Content.Find.Execute("Nonwhite", $true, $true, $false, $false, $false, $true, 1, $false, "Non-white", 1) | Out-Null
